$wb = $excel.ActiveWorkbook

# Update the "Search" worksheet data (row 2): new field label, locator type, and value
$wsSearch = $wb.Worksheets.Item("Search")
$wsSearch.Range("A2").Value = "Mensaje Contact us"
$wsSearch.Range("B2").Value = "xpath"
$wsSearch.Range("C2").Value = "//*[@class='primary-content gtm-footer-link']"

# Move the selection on the Search sheet and make it the active tab/sheet
$wsSearch.Range("C6").Select()
$wsSearch.Activate()

# Home sheet keeps its own selection at C10, but is no longer the active tab
$wsHome = $wb.Worksheets.Item("Home")
$wsHome.Range("C10").Select()

# Re-activate Search so it ends up as the selected/active sheet
$wsSearch.Activate()
